$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "21.312.44"
$ws.Range("E2").Value = "  +3.97%  "

$ws.Range("D3").Value = "1.544.74"
$ws.Range("E3").Value = "  +4.76%  "

$ws.Range("D4").Value = "'1.010"
$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").Value = "'0.9684"
$ws.Range("E5").Value = "  +0.77%  "

$ws.Range("D6").Value = "'282.32"
$ws.Range("E6").Value = "  +2.10%  "

$ws.Range("D7").Value = "'0.3631"
$ws.Range("E7").Value = "  -0.37%  "

$ws.Range("D8").Value = "'0.3210"
$ws.Range("E8").Value = "  +5.52%  "

$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").Value = "'40.91"
$ws.Range("E9").Value = "  +3.23%  "

$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "'1.114"
$ws.Range("E10").Value = "  +6.38%  "

$ws.Range("D11").Value = "'0.06949"

$ws.Range("D12").Value = "'1.005"
$ws.Range("E12").Value = "  +0.31%  "

$ws.Range("D13").Value = "'19.14"
$ws.Range("E13").Value = "  +5.71%  "

$ws.Range("D14").Value = "'5.735"
$ws.Range("E14").Value = "  +5.25%  "

$ws.Range("D15").Value = "'6.434"
$ws.Range("E15").Value = "  +4.47%  "

$ws.Range("D16").Value = "'0.00001056"
$ws.Range("E16").Value = "  +2.87%  "

$ws.Range("D17").Value = "'0.9687"
$ws.Range("E17").Value = "  +0.13%  "

$ws.Range("D18").Value = "1.546.30"
$ws.Range("E18").Value = "  +4.58%  "

$ws.Range("D19").Value = "'0.06138"
$ws.Range("E19").Value = "  +4.24%  "

$ws.Range("D20").Value = "'72.87"
$ws.Range("E20").Value = "  +5.54%  "

$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "'15.57"
$ws.Range("E21").Value = "  +7.77%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'5.754"
$ws.Range("E22").Value = "  +5.39%  "

$ws.Range("D23").Value = "'11.47"
$ws.Range("E23").Value = "  +4.73%  "

$ws.Range("D24").Value = "'2.325"
$ws.Range("E24").Value = "  +3.56%  "

$ws.Range("D25").Value = "21.383.85"
$ws.Range("E25").Value = "  +4.02%  "

$ws.Range("D26").Value = "'2.309"
$ws.Range("E26").Value = "  +8.88%  "

$ws.Range("D27").Value = "'147.44"
$ws.Range("E27").Value = "  +4.33%  "

$ws.Range("D28").Value = "'17.95"
$ws.Range("E28").Value = "  +4.58%  "

$ws.Range("D29").Value = "1.718.80"
$ws.Range("E29").Value = "  +5.15%  "

$ws.Range("D30").Value = "'119.03"
$ws.Range("E30").Value = "  +5.04%  "

$ws.Range("D31").Value = "'4.032"
$ws.Range("E31").Value = "  +3.90%  "

$ws.Range("D32").Value = "'0.8877"
$ws.Range("E32").Value = "  +9.81%  "

$ws.Range("D33").Value = "'5.266"
$ws.Range("E33").Value = "  +6.33%  "

$ws.Range("D34").Value = "'0.08097"
$ws.Range("E34").Value = "  +2.67%  "

$ws.Range("D35").Value = "'1.567"
$ws.Range("E35").Value = "  +3.96%  "

$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").Value = "'5.003"
$ws.Range("E36").Value = "  +6.01%  "

$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "'1.217"
$ws.Range("E37").Value = "  -1.74%  "

$ws.Range("D38").Value = "'0.05924"
$ws.Range("E38").Value = "  +3.02%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.02140"
$ws.Range("E39").Value = "  +5.20%  "

$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").Value = "'10.93"
$ws.Range("E40").Value = "  +5.10%  "

$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "'0.1995"
$ws.Range("E41").Value = "  +6.49%  "

$ws.Range("D42").Value = "'7.923"
$ws.Range("E42").Value = "  +2.98%  "

$ws.Range("D43").Value = "'0.9679"
$ws.Range("E43").Value = "  +0.74%  "

$ws.Range("D44").Value = "'0.5532"
$ws.Range("E44").Value = "  +5.10%  "

$ws.Range("D45").Value = "'12.69"
$ws.Range("E45").Value = "  +5.48%  "

$ws.Range("D46").Value = "'3.588"

$ws.Range("D47").Value = "'0.5518"
$ws.Range("E47").Value = "  +7.00%  "

$ws.Range("D48").Value = "'122.71"
$ws.Range("E48").Value = "  +5.00%  "

$ws.Range("D49").Value = "'1.887"
$ws.Range("E49").Value = "  +6.54%  "

$ws.Range("D50").Value = "'0.06623"
$ws.Range("E50").Value = "  +2.72%  "

$ws.Range("D51").Value = "'70.41"
$ws.Range("E51").Value = "  +5.08%  "
